$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 4
$ws.Cells.Item(4, 1).Value = '''2022-05-20'
$ws.Cells.Item(4, 2).Value = 'Tyler Penn'
$ws.Cells.Item(4, 3).Value = '[Web] [GO2bank] [ODP 2.0] Go-Live (Desktop browser)'
$ws.Cells.Item(4, 4).Value = 'https://pd.nextestate.com/browse/BUX-37261'
$ws.Cells.Item(4, 5).Value = 'M111'
$ws.Cells.Item(4, 6).Value = '''True'
$ws.Cells.Item(4, 7).Value = 'Critical'
$ws.Cells.Item(4, 8).Value = 'N/A'
$ws.Cells.Item(4, 9).Value = 'High'
$ws.Cells.Item(4, 10).Value = 'Low'
$ws.Cells.Item(4, 11).Value = 'N/A'
$ws.Cells.Item(4, 12).Value = 'Middle'
$ws.Cells.Item(4, 13).Value = 'Middle'
$ws.Cells.Item(4, 14).Value = 'N/A'
$ws.Cells.Item(4, 15).Value = '''True'
$ws.Cells.Item(4, 16).Value = '''2022-06-23'
$ws.Cells.Item(4, 17).Value = '''False'
$ws.Cells.Item(4, 18).Value = '''False'

# Row 5
$ws.Cells.Item(5, 1).Value = '''2022-05-20'
$ws.Cells.Item(5, 2).Value = 'Tyler Penn'
$ws.Cells.Item(5, 3).Value = '[Web] [GO2bank] [ODP 2.0] Go-Live (Mobile browser)'
$ws.Cells.Item(5, 4).Value = 'https://pd.nextestate.com/browse/BUX-37262'
$ws.Cells.Item(5, 5).Value = 'M111'
$ws.Cells.Item(5, 6).Value = '''True'
$ws.Cells.Item(5, 7).Value = 'Critical'
$ws.Cells.Item(5, 8).Value = 'N/A'
$ws.Cells.Item(5, 9).Value = 'High'
$ws.Cells.Item(5, 10).Value = 'Low'
$ws.Cells.Item(5, 11).Value = 'N/A'
$ws.Cells.Item(5, 12).Value = 'Middle'
$ws.Cells.Item(5, 13).Value = 'Middle'
$ws.Cells.Item(5, 14).Value = 'N/A'
$ws.Cells.Item(5, 15).Value = '''True'
$ws.Cells.Item(5, 16).Value = '''2022-06-23'
$ws.Cells.Item(5, 17).Value = '''False'
$ws.Cells.Item(5, 18).Value = '''False'

# Row 6
$ws.Cells.Item(6, 1).Value = '''2022-05-20'
$ws.Cells.Item(6, 2).Value = 'Tyler Penn'
$ws.Cells.Item(6, 3).Value = '[iOS/Android][GO2bank] [ODP 2.0] Go-Live'
$ws.Cells.Item(6, 4).Value = 'https://pd.nextestate.com/browse/BMAPP-21834'
$ws.Cells.Item(6, 5).Value = 'M111'
$ws.Cells.Item(6, 6).Value = '''True'
$ws.Cells.Item(6, 7).Value = 'Critical'
$ws.Cells.Item(6, 8).Value = 'N/A'
$ws.Cells.Item(6, 9).Value = 'High'
$ws.Cells.Item(6, 10).Value = 'Low'
$ws.Cells.Item(6, 11).Value = 'N/A'
$ws.Cells.Item(6, 12).Value = 'Middle'
$ws.Cells.Item(6, 13).Value = 'Middle'
$ws.Cells.Item(6, 14).Value = 'N/A'
$ws.Cells.Item(6, 15).Value = '''True'
$ws.Cells.Item(6, 16).Value = '''2022-06-23'
$ws.Cells.Item(6, 17).Value = '''False'
$ws.Cells.Item(6, 18).Value = '''False'

# Row 7
$ws.Cells.Item(7, 1).Value = '''2022-06-02'
$ws.Cells.Item(7, 2).Value = 'Tyler Penn'
$ws.Cells.Item(7, 3).Value = 'ODP 2.0 : [TRIGGER FOR NTKey105] Tier Reinstated'
$ws.Cells.Item(7, 4).Value = 'https://pd.nextestate.com/browse/GBOS-62115'
$ws.Cells.Item(7, 5).Value = 'M111'
$ws.Cells.Item(7, 6).Value = '''True'
$ws.Cells.Item(7, 7).Value = 'High'
$ws.Cells.Item(7, 8).Value = 'N/A'
$ws.Cells.Item(7, 9).Value = 'High'
$ws.Cells.Item(7, 10).Value = 'N/A'
$ws.Cells.Item(7, 11).Value = 'N/A'
$ws.Cells.Item(7, 12).Value = 'Low'
$ws.Cells.Item(7, 13).Value = 'N/A'
$ws.Cells.Item(7, 14).Value = 'High'
$ws.Cells.Item(7, 15).Value = '''True'
$ws.Cells.Item(7, 16).Value = '''2022-06-23'
$ws.Cells.Item(7, 17).Value = '''False'
$ws.Cells.Item(7, 18).Value = '''False'

# Row 8
$ws.Cells.Item(8, 1).Value = '''2022-06-07'
$ws.Cells.Item(8, 2).Value = 'Smitha Jonnala'
$ws.Cells.Item(8, 3).Value = '[WMMC] OAuth redirect is not working as expected'
$ws.Cells.Item(8, 4).Value = 'https://pd.nextestate.com/browse/GB-80048'
$ws.Cells.Item(8, 5).Value = 'M111'
$ws.Cells.Item(8, 6).Value = '''True'
$ws.Cells.Item(8, 7).Value = 'N/A'
$ws.Cells.Item(8, 8).Value = 'N/A'
$ws.Cells.Item(8, 9).Value = 'Middle'
$ws.Cells.Item(8, 10).Value = 'N/A'
$ws.Cells.Item(8, 11).Value = 'N/A'
$ws.Cells.Item(8, 12).Value = 'N/A'
$ws.Cells.Item(8, 13).Value = 'N/A'
$ws.Cells.Item(8, 14).Value = 'Middle'
$ws.Cells.Item(8, 15).Value = '''True'
$ws.Cells.Item(8, 16).Value = '''2022-06-23'
$ws.Cells.Item(8, 17).Value = '''False'
$ws.Cells.Item(8, 18).Value = '''False'

# Row 9
$ws.Cells.Item(9, 1).Value = '''2022-06-15'
$ws.Cells.Item(9, 2).Value = 'Mark Sun'
$ws.Cells.Item(9, 3).Value = '[SCC] Block non-SSN DDA applying for an SCC'
$ws.Cells.Item(9, 4).Value = 'https://pd.nextestate.com/browse/GBOS-62345'
$ws.Cells.Item(9, 5).Value = 'M111'
$ws.Cells.Item(9, 6).Value = '''False'
$ws.Cells.Item(9, 7).Value = 'Critical'
$ws.Cells.Item(9, 8).Value = 'Low'
$ws.Cells.Item(9, 9).Value = 'Middle'
$ws.Cells.Item(9, 10).Value = 'N/A'
$ws.Cells.Item(9, 11).Value = 'N/A'
$ws.Cells.Item(9, 12).Value = 'High'
$ws.Cells.Item(9, 13).Value = 'High'
$ws.Cells.Item(9, 14).Value = 'High'
$ws.Cells.Item(9, 15).Value = '''False'
$ws.Cells.Item(9, 16).Value = '''2022-06-23'
$ws.Cells.Item(9, 17).Value = '''False'
$ws.Cells.Item(9, 18).Value = '''False'

# Row 10
$ws.Cells.Item(10, 1).Value = '''2022-06-22'
$ws.Cells.Item(10, 2).Value = 'Tyler Penn'
$ws.Cells.Item(10, 3).Value = '[Gateway] [GO2bank] [app upgrade] Set app force upgrade'
$ws.Cells.Item(10, 4).Value = 'https://pd.nextestate.com/browse/BUX-38827'
$ws.Cells.Item(10, 5).Value = 'M111'
$ws.Cells.Item(10, 6).Value = '''False'
$ws.Cells.Item(10, 7).Value = 'Critical'
$ws.Cells.Item(10, 8).Value = 'Middle'
$ws.Cells.Item(10, 9).Value = 'High'
$ws.Cells.Item(10, 10).Value = 'Middle'
$ws.Cells.Item(10, 11).Value = 'N/A'
$ws.Cells.Item(10, 12).Value = 'N/A'
$ws.Cells.Item(10, 13).Value = 'Middle'
$ws.Cells.Item(10, 14).Value = 'N/A'
$ws.Cells.Item(10, 15).Value = '''True'
$ws.Cells.Item(10, 16).Value = '''2022-06-23'
$ws.Cells.Item(10, 17).Value = '''False'
$ws.Cells.Item(10, 18).Value = '''False'

# Row 11
$ws.Cells.Item(11, 1).Value = '''2022-05-23'
$ws.Cells.Item(11, 2).Value = 'Megan Ackling'
$ws.Cells.Item(11, 3).Value = '[QA Testing]ATM and AFT Limit Decrease for Regular Season - CoID 164'
$ws.Cells.Item(11, 4).Value = 'https://pd.nextestate.com/browse/GBOS-62061'
$ws.Cells.Item(11, 5).Value = 'M111'
$ws.Cells.Item(11, 6).Value = '''False'
$ws.Cells.Item(11, 7).Value = 'High'
$ws.Cells.Item(11, 8).Value = 'N/A'
$ws.Cells.Item(11, 9).Value = 'N/A'
$ws.Cells.Item(11, 10).Value = 'N/A'
$ws.Cells.Item(11, 11).Value = 'N/A'
$ws.Cells.Item(11, 12).Value = 'Middle'
$ws.Cells.Item(11, 13).Value = 'N/A'
$ws.Cells.Item(11, 14).Value = 'N/A'
$ws.Cells.Item(11, 15).Value = '''False'
$ws.Cells.Item(11, 16).Value = '''2022-06-23'
$ws.Cells.Item(11, 17).Value = '''False'
$ws.Cells.Item(11, 18).Value = '''False'

# Row 12
$ws.Cells.Item(12, 1).Value = '''2022-06-14'
$ws.Cells.Item(12, 2).Value = 'Sarath Krishnan'
$ws.Cells.Item(12, 3).Value = '[GFT]Update the Utility Core API to get Bin info to add Program code'
$ws.Cells.Item(12, 4).Value = 'https://pd.nextestate.com/browse/GBOS-62197'
$ws.Cells.Item(12, 5).Value = 'M111'
$ws.Cells.Item(12, 6).Value = '''False'
$ws.Cells.Item(12, 7).Value = 'N/A'
$ws.Cells.Item(12, 8).Value = 'High'
$ws.Cells.Item(12, 9).Value = 'High'
$ws.Cells.Item(12, 10).Value = 'N/A'
$ws.Cells.Item(12, 11).Value = 'N/A'
$ws.Cells.Item(12, 12).Value = 'N/A'
$ws.Cells.Item(12, 13).Value = 'High'
$ws.Cells.Item(12, 14).Value = 'N/A'
$ws.Cells.Item(12, 15).Value = '''False'
$ws.Cells.Item(12, 16).Value = '''2022-06-23'
$ws.Cells.Item(12, 17).Value = '''False'
$ws.Cells.Item(12, 18).Value = '''False'

# Row 13
$ws.Cells.Item(13, 1).Value = '''2022-05-20'
$ws.Cells.Item(13, 2).Value = 'Megan Ackling'
$ws.Cells.Item(13, 3).Value = 'Walgreens and CVS FeeSplit Setup for eCash Setup - GBOS'
$ws.Cells.Item(13, 4).Value = 'https://pd.nextestate.com/browse/FEAS-24667'
$ws.Cells.Item(13, 5).Value = 'M111'
$ws.Cells.Item(13, 6).Value = '''False'
$ws.Cells.Item(13, 7).Value = 'N/A'
$ws.Cells.Item(13, 8).Value = 'High'
$ws.Cells.Item(13, 9).Value = 'High'
$ws.Cells.Item(13, 10).Value = 'N/A'
$ws.Cells.Item(13, 11).Value = 'N/A'
$ws.Cells.Item(13, 12).Value = 'N/A'
$ws.Cells.Item(13, 13).Value = 'N/A'
$ws.Cells.Item(13, 14).Value = 'N/A'
$ws.Cells.Item(13, 15).Value = '''False'
$ws.Cells.Item(13, 16).Value = '''2022-06-23'
$ws.Cells.Item(13, 17).Value = '''False'
$ws.Cells.Item(13, 18).Value = '''False'

# Row 14
$ws.Cells.Item(14, 1).Value = '''2022-06-06'
$ws.Cells.Item(14, 2).Value = 'Sarath Krishnan'
$ws.Cells.Item(14, 3).Value = '[GFT][MCsend]Create Adjusmentypes for Partner A2A and P2P'
$ws.Cells.Item(14, 4).Value = 'https://pd.nextestate.com/browse/GBOS-61750'
$ws.Cells.Item(14, 5).Value = 'M111'
$ws.Cells.Item(14, 6).Value = '''False'
$ws.Cells.Item(14, 7).Value = 'N/A'
$ws.Cells.Item(14, 8).Value = 'High'
$ws.Cells.Item(14, 9).Value = 'High'
$ws.Cells.Item(14, 10).Value = 'N/A'
$ws.Cells.Item(14, 11).Value = 'N/A'
$ws.Cells.Item(14, 12).Value = 'N/A'
$ws.Cells.Item(14, 13).Value = 'N/A'
$ws.Cells.Item(14, 14).Value = 'N/A'
$ws.Cells.Item(14, 15).Value = '''False'
$ws.Cells.Item(14, 16).Value = '''2022-06-23'
$ws.Cells.Item(14, 17).Value = '''False'
$ws.Cells.Item(14, 18).Value = '''False'

# Row 15
$ws.Cells.Item(15, 1).Value = '''2022-06-16'
$ws.Cells.Item(15, 2).Value = 'Dennis Wiles'
$ws.Cells.Item(15, 3).Value = '[GO2bank][Web] Add Domains to CSP Whitelisting'
$ws.Cells.Item(15, 4).Value = 'https://pd.nextestate.com/browse/BUX-38699'
$ws.Cells.Item(15, 5).Value = 'M111'
$ws.Cells.Item(15, 6).Value = '''False'
$ws.Cells.Item(15, 7).Value = 'N/A'
$ws.Cells.Item(15, 8).Value = 'N/A'
$ws.Cells.Item(15, 9).Value = 'High'
$ws.Cells.Item(15, 10).Value = 'High'
$ws.Cells.Item(15, 11).Value = 'N/A'
$ws.Cells.Item(15, 12).Value = 'N/A'
$ws.Cells.Item(15, 13).Value = 'N/A'
$ws.Cells.Item(15, 14).Value = 'N/A'
$ws.Cells.Item(15, 15).Value = '''False'
$ws.Cells.Item(15, 16).Value = '''2022-06-23'
$ws.Cells.Item(15, 17).Value = '''False'
$ws.Cells.Item(15, 18).Value = '''False'

# Row 16
$ws.Cells.Item(16, 1).Value = '''2022-06-13'
$ws.Cells.Item(16, 2).Value = 'Rita Webb'
$ws.Cells.Item(16, 3).Value = '[COFO]  Research Extensibility of Current Mobile 2FA Service'
$ws.Cells.Item(16, 4).Value = '[GBOS-61873] [COFO][SPIKE] Research Extensibility of Current Mobile 2FA Service - GDCJira (nextestate.com)'
$ws.Cells.Item(16, 5).Value = 'M111'
$ws.Cells.Item(16, 6).Value = '''False'
$ws.Cells.Item(16, 7).Value = 'N/A'
$ws.Cells.Item(16, 8).Value = 'Middle'
$ws.Cells.Item(16, 9).Value = 'High'
$ws.Cells.Item(16, 10).Value = 'N/A'
$ws.Cells.Item(16, 11).Value = 'N/A'
$ws.Cells.Item(16, 12).Value = 'Middle'
$ws.Cells.Item(16, 13).Value = 'Low'
$ws.Cells.Item(16, 14).Value = 'Low'
$ws.Cells.Item(16, 15).Value = '''False'
$ws.Cells.Item(16, 16).Value = '''2022-06-23'
$ws.Cells.Item(16, 17).Value = '''False'
$ws.Cells.Item(16, 18).Value = '''False'

# Row 17
$ws.Cells.Item(17, 1).Value = '''2022-06-16'
$ws.Cells.Item(17, 2).Value = 'Mark Sun'
$ws.Cells.Item(17, 3).Value = '[QA ONLY] SCC declines can re-apply after 30 days'
$ws.Cells.Item(17, 4).Value = 'https://pd.nextestate.com/browse/BMAPP-22159'
$ws.Cells.Item(17, 5).Value = 'M111'
$ws.Cells.Item(17, 6).Value = '''False'
$ws.Cells.Item(17, 7).Value = 'N/A'
$ws.Cells.Item(17, 8).Value = 'Low'
$ws.Cells.Item(17, 9).Value = 'High'
$ws.Cells.Item(17, 10).Value = 'N/A'
$ws.Cells.Item(17, 11).Value = 'N/A'
$ws.Cells.Item(17, 12).Value = 'Low'
$ws.Cells.Item(17, 13).Value = 'Low'
$ws.Cells.Item(17, 14).Value = 'N/A'
$ws.Cells.Item(17, 15).Value = '''False'
$ws.Cells.Item(17, 16).Value = '''2022-06-23'
$ws.Cells.Item(17, 17).Value = '''False'
$ws.Cells.Item(17, 18).Value = '''False'

# Row 18
$ws.Cells.Item(18, 1).Value = '''2022-06-16'
$ws.Cells.Item(18, 2).Value = 'Mark Sun'
$ws.Cells.Item(18, 3).Value = '[QA only] SCC declines can re-apply after 30 days'
$ws.Cells.Item(18, 4).Value = 'https://pd.nextestate.com/browse/BUX-38697'
$ws.Cells.Item(18, 5).Value = 'M111'
$ws.Cells.Item(18, 6).Value = '''False'
$ws.Cells.Item(18, 7).Value = 'N/A'
$ws.Cells.Item(18, 8).Value = 'Low'
$ws.Cells.Item(18, 9).Value = 'High'
$ws.Cells.Item(18, 10).Value = 'N/A'
$ws.Cells.Item(18, 11).Value = 'N/A'
$ws.Cells.Item(18, 12).Value = 'Low'
$ws.Cells.Item(18, 13).Value = 'Low'
$ws.Cells.Item(18, 14).Value = 'N/A'
$ws.Cells.Item(18, 15).Value = '''False'
$ws.Cells.Item(18, 16).Value = '''2022-06-23'
$ws.Cells.Item(18, 17).Value = '''False'
$ws.Cells.Item(18, 18).Value = '''False'
